$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H2 from a plain number to a text value "86(12)"
$ws.Range("H2").Value = "86(12)"

# Move the active selection to H2
$ws.Range("H2").Select()
